$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.996.65"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.828.01"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'312.17"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "'0.4619"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.3703"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").Value = "'0.07339"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "'0.8758"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").Value = "'0.07925"
$ws.Range("E11").Value = "  +4.52%  "
$ws.Range("D12").Value = "'19.83"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "1.781.07"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "'5.343"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "'6.561"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'91.42"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "'0.000008874"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'14.82"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "27.279.56"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D24").Value = "2.057.88"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").Value = "'153.04"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'1.850"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").Value = "'18.44"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("D29").Value = "'5.153"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").Value = "'115.52"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "'0.08894"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "'2.963"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'0.7301"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'4.436"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "'1.129"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Value = "'2.472"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("D38").Value = "'1.069"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'0.05234"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "'2.947"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "'7.118"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'0.5164"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'8.187"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").Value = "'0.4852"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "'1.005"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "'10.17"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'102.68"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'1.635"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'0.06200"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "'65.00"
$ws.Range("E51").Value = "  +0.67%  "
